$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Sheet1 "Forecast Comparison" updates ---
$ws1.Range("D2").Value = 9
$ws1.Range("H2").Value = 19.02
$ws1.Range("L2").Value = 1.14
$ws1.Range("D3").Value = 10
$ws1.Range("H3").Value = 16.25
$ws1.Range("L3").Value = 1.15
$ws1.Range("D4").Value = 13
$ws1.Range("H4").Value = 11.79
$ws1.Range("L4").Value = 1.12
$ws1.Range("D5").Value = 15
$ws1.Range("H5").Value = 9.25
$ws1.Range("L5").Value = 0.96
$ws1.Range("D6").Value = 15
$ws1.Range("H6").Value = 8.470000000000001
$ws1.Range("L6").Value = 0.88
$ws1.Range("D7").Value = 13
$ws1.Range("H7").Value = 8.619999999999999
$ws1.Range("L7").Value = 1.19
$ws1.Range("D8").Value = 12
$ws1.Range("H8").Value = 8.109999999999999
$ws1.Range("L8").Value = 0.95
$ws1.Range("D9").Value = 14
$ws1.Range("H9").Value = 6.2
$ws1.Range("L9").Value = 1.06
$ws1.Range("D10").Value = 16
$ws1.Range("H10").Value = 4.49
$ws1.Range("L10").Value = 0.84
$ws1.Range("D11").Value = 16
$ws1.Range("H11").Value = 3.45
$ws1.Range("L11").Value = 1.15
$ws1.Range("D12").Value = 14
$ws1.Range("H12").Value = 2.75
$ws1.Range("L12").Value = 0.88
$ws1.Range("D13").Value = 13
$ws1.Range("H13").Value = 1.88
$ws1.Range("L13").Value = 0.95
$ws1.Range("D14").Value = 15
$ws1.Range("H14").Value = 0.8
$ws1.Range("J14").Value = "Urgent"
$ws1.Range("L14").Value = 1.13
$ws1.Range("D15").Value = 17
$ws1.Range("H15").Value = 0
$ws1.Range("I15").Value = "High"
$ws1.Range("J15").Value = "Urgent"
$ws1.Range("L15").Value = 0.88
$ws1.Range("D16").Value = 17
$ws1.Range("H16").Value = 0
$ws1.Range("I16").Value = "High"
$ws1.Range("J16").Value = "Urgent"
$ws1.Range("L16").Value = 0.86
$ws1.Range("D17").Value = 16
$ws1.Range("H17").Value = 0
$ws1.Range("I17").Value = "High"
$ws1.Range("J17").Value = "Urgent"
$ws1.Range("L17").Value = 1.13

# --- Sheet2 "Summary" updates (stored as text strings) ---
$ws2.Range("B9").Value = "229"
$ws2.Range("B10").Value = "102"
$ws2.Range("B11").Value = "48"
$ws2.Range("B12").Value = "18"
$ws2.Range("B14").Value = "9"

Write-Host "Edit applied successfully"
